$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely
$meta.Range("A11").EntireRow.Delete()

# --- Sheet "Elements" (sheet2) ---
$elem = $wb.Worksheets.Item("Elements")

# Root element row (row 2): Short/Definition get profile-specific text instead of generic extension text
$elem.Range("K2").Value = "Enrollment PCP Zip On Enrollment"
$elem.Range("L2").Value = "Original primary care physician (PCP) zip code as reporting on the eligibility record"
